# Auto-generated: update Price (D) and Volume(1h) (E) columns per crypto data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.968.39'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '2.623.70'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.42'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.38'
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('E9').Value = '  -1.41%  '
$ws.Range('E10').Value = '  -1.51%  '
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.59'
$ws.Range('E14').Value = '  -0.47%  '
$ws.Range('D15').Value = '3.109.56'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '66.918.25'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '2.624.76'
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.15'
$ws.Range('E19').Value = '  +4.81%  '
$ws.Range('E20').Value = '  +7.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '355.82'
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('E23').Value = '  -2.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.89'
$ws.Range('E24').Value = '  +8.45%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.91'
$ws.Range('E26').Value = '  -5.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '70.05'
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('D28').Value = '2.761.67'
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000100'
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '547.79'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.90'
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('E33').Value = '  -2.03%  '
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('E35').Value = '  +5.02%  '
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('E37').Value = '  -4.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '156.48'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.02'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.366'
$ws.Range('E40').Value = '  -1.45%  '
$ws.Range('E41').Value = '  -2.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.15'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.40'
$ws.Range('E46').Value = '  -4.85%  '
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.577'
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '151.22'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('E51').Value = '  -0.61%  '
